$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.293070316314697
$ws.Range("B1").Value = 2.044524431228638
$ws.Range("C1").Value = 5.447574615478516
$ws.Range("D1").Value = 1.911295056343079
$ws.Range("E1").Value = 1.098896980285645
